$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("DS1")
$ws2 = $wb.Worksheets.Item("DS2")

# Rename the "kVec" header column to "FC" on both data sheets.
$ws1.Range("B1").Value = "FC"
$ws2.Range("B1").Value = "FC"

# Move the active selection / active sheet: DS1 becomes the active tab with
# B1 selected, while DS2 (previously active) keeps a new selection at G2.
$ws2.Range("G2").Select()
$ws1.Activate()
$ws1.Range("B1").Select()
